# Applies "region type, pressure loss, report correction" edits.
# Each data value lives in column 4 of a small 4-column key/value table;
# we target cells by (table index, row index) so that duplicate values
# (e.g. "305" appearing in two different tables) are updated correctly
# and independently.

$d = $word.ActiveDocument
$tables = $d.Tables

function Set-CellText($tableIndex, $rowIndex, $colIndex, $newText) {
    $cell = $tables.Item($tableIndex).Cell($rowIndex, $colIndex)
    $cell.Range.Text = $newText
}

# Header table: Enquiry/Date/Model
Set-CellText 1 2 4 "11/25/2019, 10:22 AM"   # Date
Set-CellText 1 3 4 "TAC S2 C3"              # Model

# Capacity table
Set-CellText 2 2 4 "114"                    # Capacity (TR)

# Chilled water circuit table
Set-CellText 3 2 4 "68.8"                   # Chilled water flow
Set-CellText 3 6 4 "1.2"                    # Chilled water circuit pressure loss
Set-CellText 3 7 4 "125"                    # Chilled water Connection diameter

# Cooling water circuit table
Set-CellText 4 2 4 "114"                    # Cooling water flow
Set-CellText 4 4 4 "37.1"                   # Cooling water outlet temperature
Set-CellText 4 7 4 "2.2"                    # Cooling water circuit pressure loss
Set-CellText 4 8 4 "150"                    # Cooling water Connection diameter

# Steam circuit table
Set-CellText 5 3 4 "400.5"                  # Steam Consumption(+/-3%)
Set-CellText 5 6 4 "65"                     # Connection - Inlet diameter

# Electrical data table
Set-CellText 6 3 4 "7.6"                    # Power consumption
Set-CellText 6 4 4 "2.2( 6 )"               # Absorbent pump rating

# Physical data table
Set-CellText 7 2 4 "3140"                   # Length
Set-CellText 7 3 4 "2140"                   # Width
Set-CellText 7 4 4 "2750"                   # Height
Set-CellText 7 5 4 "6.8"                    # Operating weight
Set-CellText 7 6 4 "6.4"                    # Shipping weight
Set-CellText 7 7 4 "9.4"                    # Flooded weight
Set-CellText 7 8 4 "5.5"                    # Dry weight
Set-CellText 7 9 4 "2560"                   # Tube cleaning space
